$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '62.824.18'
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.969.62'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +1.39%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '593.96'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.20%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '145.95'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('E7').Value = '  -0.04%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '2.967.51'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +1.35%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.506'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +0.41%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '7.23'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +3.22%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.145'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +1.88%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.445'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +1.04%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.0000237'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +5.40%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '33.29'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.23%  '
$ws.Range('E15').Value = '  -0.34%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.459.82'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +1.30%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '62.738.63'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +2.67%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '6.72'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.16%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '2.953.27'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.81%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '442.54'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +2.03%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '13.46'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E22').Value = '  -1.06%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '7.08'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.30%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '11.32'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +2.13%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '81.82'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +0.42%  '
$ws.Range('E26').Value = '  -3.11%  '
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('E28').Value = '  +0.03%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '7.19'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +3.35%  '
$ws.Range('E30').Value = '  +0.01%  '
$ws.Range('E31').Value = '  -5.19%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.0₃0936'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +7.89%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.110'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '26.71'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.24%  '
$ws.Range('E35').Value = '  +0.02%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -1.56%  '
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('E38').Value = '  +1.42%  '
$ws.Range('E39').Value = '  +3.15%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '49.51'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.84%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '8.57'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.35%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.118'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -4.27%  '
$ws.Range('E43').Value = '  -0.93%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '39.89'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -5.08%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '2.740.01'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +1.17%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '136.12'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +1.59%  '
$ws.Range('E47').Value = '  -1.57%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '364.02'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('E50').Value = '  -0.14%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '23.05'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -3.12%  '
